# Update "想去人数" (F column) counts on the 展览, 演出, and 全部类型 sheets
# to reflect the latest scraped totals (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 608
$ws1.Range("F5").Value = 649
$ws1.Range("F10").Value = 1144
$ws1.Range("F11").Value = 586
$ws1.Range("F12").Value = 348
$ws1.Range("F15").Value = 303
$ws1.Range("F20").Value = 533
$ws1.Range("F22").Value = 508

# 演出 sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 637

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 73
$ws4.Range("F7").Value = 608
$ws4.Range("F9").Value = 649
$ws4.Range("F14").Value = 1144
$ws4.Range("F15").Value = 586
$ws4.Range("F18").Value = 348
$ws4.Range("F20").Value = 637
$ws4.Range("F23").Value = 303
$ws4.Range("F32").Value = 533
$ws4.Range("F34").Value = 508
